# Update "Alert Parameters (working).xlsx" - Third Iteration sheet
# Rename the "Lee3" building identifier to "Lee_III" in the bldg_std (F)
# and database (H) columns for the relevant alert rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Third Iteration")

# F column: bldg_std values "LEE3" -> "LEE_III"
$ws.Range("F24").Value = "LEE_III"
$ws.Range("F25").Value = "LEE_III"
$ws.Range("F26").Value = "LEE_III"
$ws.Range("F31").Value = "LEE_III"
$ws.Range("F32").Value = "LEE_III"
$ws.Range("F33").Value = "LEE_III"
$ws.Range("F34").Value = "LEE_III"

# H column: database values "CEVAC_LEE3_..." -> "CEVAC_LEE_III_..."
$ws.Range("H24").Value = "CEVAC_LEE_III_TEMP_LATEST"
$ws.Range("H25").Value = "CEVAC_LEE_III_POWER_LATEST"
$ws.Range("H26").Value = "CEVAC_LEE_III_IAQ_LATEST"
$ws.Range("H31").Value = "CEVAC_LEE_III_TEMP_DAY"
$ws.Range("H32").Value = "CEVAC_LEE_III_TEMP_LATEST"
$ws.Range("H33").Value = "CEVAC_LEE_III_TEMP_DAY"
$ws.Range("H34").Value = "CEVAC_LEE_III_TEMP_LATEST"

# Update the sheet view: scroll so row 10 is at the top, and select H14
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H14").Select()
